# Update the "Demographics and Preoperative Characteristics" table values
# for the "uninsured" cohort. Each (old, new) pair below corresponds to a
# single, uniquely-occurring cell value in the document, so a literal
# (non-wildcard) Find/Replace-All is unambiguous for every pair.
#
# NOTE: pairs are applied in the same top-to-bottom order they appear in
# the document. This matters for one pair: "26.0 (0.06)" -> "25 (0.06)"
# produces text that coincidentally matches another row's *old* value
# ("25 (0.06)" -> "20 (0.05)"). Because that other row appears earlier in
# the document and is processed first, its text is already changed to
# "20 (0.05)" by the time we create the new "25 (0.06)" string, so no
# double-replacement occurs.
$d = $word.ActiveDocument

$replacements = @(
    @("34.16 +/- 0.07", "36.19 +/- 0.06"),
    @("484 (1.10)", "0 (0.00)"),
    @("1,979 (4.49)", "1,765 (4.52)"),
    @("678 (1.54)", "580 (1.49)"),
    @("17,001 (38.57)", "15,036 (38.53)"),
    @("26,859 (60.93)", "23,818 (61.03)"),
    @("222 (0.50)", "172 (0.44)"),
    @("1,107 (2.51)", "957 (2.45)"),
    @("4,391 (9.96)", "4,078 (10.45)"),
    @("15,138 (34.34)", "13,143 (33.68)"),
    @("239 (0.54)", "205 (0.53)"),
    @("2,397 (5.44)", "2,088 (5.35)"),
    @("2,240 (5.08)", "1,869 (4.79)"),
    @("18,570 (42.13)", "16,686 (42.76)"),
    @("44,082 (100.00)", "39,026 (100.00)"),
    @("5,108 (11.59)", "4,607 (11.80)"),
    @("16,195 (36.74)", "14,826 (37.99)"),
    @("22,779 (51.67)", "19,593 (50.20)"),
    @("6,413 (14.55)", "5,658 (14.50)"),
    @("5,932 (13.46)", "5,178 (13.27)"),
    @("23,044 (52.28)", "20,801 (53.30)"),
    @("8,693 (19.72)", "7,389 (18.93)"),
    @("16,117 (36.56)", "14,498 (37.15)"),
    @("12,094 (27.44)", "10,706 (27.43)"),
    @("9,908 (22.48)", "8,683 (22.25)"),
    @("5,963 (13.53)", "5,139 (13.17)"),
    @("204 (0.46)", "169 (0.43)"),
    @("6,879 (15.61)", "5,791 (14.84)"),
    @("25 (0.06)", "20 (0.05)"),
    @("9,960 (22.59)", "8,958 (22.95)"),
    @("40.0 (0.09)", "40 (0.10)"),
    @("718.0 (1.63)", "707 (1.81)"),
    @("110.0 (0.25)", "101 (0.26)"),
    @("11.0 (0.02)", "9 (0.02)"),
    @("26.0 (0.06)", "25 (0.06)"),
    @("62.0 (0.14)", "52 (0.13)"),
    @("896.0 (2.03)", "868 (2.22)"),
    @("1,728.0 (3.92)", "1,637 (4.19)"),
    @("212.0 (0.48)", "195 (0.50)"),
    @("4,250.0 (9.64)", "4,029 (10.32)"),
    @("3,908.0 (8.87)", "3,711 (9.51)"),
    @("1,850.0 (4.20)", "1,643 (4.21)"),
    @("2,959.0 (6.71)", "2,846 (7.29)"),
    @("102.0 (0.23)", "86 (0.22)"),
    @("546.0 (1.24)", "514 (1.32)"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    # MatchCase=$true, MatchWholeWord=$false, MatchWildcards=$false,
    # MatchSoundsLike=$false, MatchAllWordForms=$false, Forward=$true,
    # Wrap=1 (wdFindContinue), Format=$false, Replace=2 (wdReplaceAll)
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
